$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp
$ws.Range("A1").Value2 = "Datos actualizados a 22 de Septiembre de 2020 a las 16:53"

# Update per-country statistics (new data pulled for this refresh)
    # Estados Unidos (row 4)
    $ws.Cells.Item(4,2).Value2 = 7049165
    $ws.Cells.Item(4,3).Value2 = 2949
    $ws.Cells.Item(4,4).Value2 = 4301089
    $ws.Cells.Item(4,5).Value2 = 2543443
    $ws.Cells.Item(4,7).Value2 = 127
    $ws.Cells.Item(4,8).Value2 = 204633
    # India (row 5)
    $ws.Cells.Item(5,2).Value2 = 5580286
    $ws.Cells.Item(5,3).Value2 = 20181
    $ws.Cells.Item(5,4).Value2 = 4509924
    $ws.Cells.Item(5,5).Value2 = 981245
    $ws.Cells.Item(5,7).Value2 = 152
    $ws.Cells.Item(5,8).Value2 = 89117
    # Chile (row 15)
    $ws.Cells.Item(15,2).Value2 = 448523
    $ws.Cells.Item(15,3).Value2 = 1055
    $ws.Cells.Item(15,4).Value2 = 423176
    $ws.Cells.Item(15,5).Value2 = 13026
    $ws.Cells.Item(15,7).Value2 = 23
    $ws.Cells.Item(15,8).Value2 = 12321
    # Irak (row 20)
    $ws.Cells.Item(20,2).Value2 = 327580
    $ws.Cells.Item(20,3).Value2 = 4724
    $ws.Cells.Item(20,4).Value2 = 261757
    $ws.Cells.Item(20,5).Value2 = 57141
    $ws.Cells.Item(20,7).Value2 = 57
    $ws.Cells.Item(20,8).Value2 = 8682
    # Alemania (row 25)
    $ws.Cells.Item(25,2).Value2 = 276226
    $ws.Cells.Item(25,3).Value2 = 675
    $ws.Cells.Item(25,5).Value2 = 20441
    # Portugal (row 52)
    $ws.Cells.Item(52,2).Value2 = 69663
    $ws.Cells.Item(52,3).Value2 = 463
    $ws.Cells.Item(52,4).Value2 = 45974
    $ws.Cells.Item(52,5).Value2 = 21764
    $ws.Cells.Item(52,7).Value2 = 5
    $ws.Cells.Item(52,8).Value2 = 1925
    # Singapur (row 57)
    $ws.Cells.Item(57,4).Value2 = 57262
    $ws.Cells.Item(57,5).Value2 = 338
    # Republica de Macedonia (row 86)
    $ws.Cells.Item(86,2).Value2 = 16867
    $ws.Cells.Item(86,3).Value2 = 87
    $ws.Cells.Item(86,4).Value2 = 14084
    $ws.Cells.Item(86,5).Value2 = 2078
    $ws.Cells.Item(86,7).Value2 = 5
    $ws.Cells.Item(86,8).Value2 = 705
    # Madagascar (row 87)
    $ws.Cells.Item(87,2).Value2 = 16136
    $ws.Cells.Item(87,3).Value2 = 63
    $ws.Cells.Item(87,4).Value2 = 14743
    $ws.Cells.Item(87,5).Value2 = 1167
    $ws.Cells.Item(87,7).Value2 = 1
    $ws.Cells.Item(87,8).Value2 = 226
    # Zambia (row 91)
    $ws.Cells.Item(91,2).Value2 = 14389
    $ws.Cells.Item(91,3).Value2 = 214
    $ws.Cells.Item(91,5).Value2 = 429
    # Noruega (row 93)
    $ws.Cells.Item(93,2).Value2 = 13075
    $ws.Cells.Item(93,3).Value2 = 70
    $ws.Cells.Item(93,5).Value2 = 2437
    # Albania (row 94)
    $ws.Cells.Item(94,2).Value2 = 12656
    $ws.Cells.Item(94,3).Value2 = 121
    $ws.Cells.Item(94,5).Value2 = 5294
    $ws.Cells.Item(94,7).Value2 = 3
    $ws.Cells.Item(94,8).Value2 = 367
    # Namibia (row 96)
    $ws.Cells.Item(96,2).Value2 = 10607
    $ws.Cells.Item(96,3).Value2 = 81
    $ws.Cells.Item(96,4).Value2 = 8359
    $ws.Cells.Item(96,5).Value2 = 2132
    $ws.Cells.Item(96,7).Value2 = 3
    $ws.Cells.Item(96,8).Value2 = 116
    # Birmania (row 113)
    $ws.Cells.Item(113,2).Value2 = 6743
    $ws.Cells.Item(113,3).Value2 = 592
    $ws.Cells.Item(113,4).Value2 = 1951
    $ws.Cells.Item(113,5).Value2 = 4677
    $ws.Cells.Item(113,7).Value2 = 17
    $ws.Cells.Item(113,8).Value2 = 115
    # Cuba (row 119)
    $ws.Cells.Item(119,2).Value2 = 5222
    $ws.Cells.Item(119,3).Value2 = 81
    $ws.Cells.Item(119,4).Value2 = 4506
    $ws.Cells.Item(119,5).Value2 = 599
    $ws.Cells.Item(119,7).Value2 = 1
    $ws.Cells.Item(119,8).Value2 = 117
    # Sri Lanka (row 142)
    $ws.Cells.Item(142,2).Value2 = 3312
    $ws.Cells.Item(142,3).Value2 = 13
    $ws.Cells.Item(142,5).Value2 = 181
    # Mauricio (row 180)
    $ws.Cells.Item(180,2).Value2 = 367
    $ws.Cells.Item(180,3).Value2 = 1
    $ws.Cells.Item(180,4).Value2 = 339
    # Bonaire, San Eustaquio y Saba (row 201)
    $ws.Cells.Item(201,2).Value2 = 54
    $ws.Cells.Item(201,3).Value2 = 18
    $ws.Cells.Item(201,4).Value2 = 18
    $ws.Cells.Item(201,5).Value2 = 35

# Re-sort the data range by "Casos totales" (column B) descending,
# same as the published sheet keeps it ranked.
$dataRange = $ws.Range("A4:H219")
$sortKey = $ws.Range("B4")
$dataRange.Sort($sortKey, 2)
